# Add the two DADA2 rows (DADA2_Species, DADA2_Taxonomy) into the Mean_F1_table
# sheet, keeping the existing alphabetical ordering of the "Type" column.
# This pushes the previous rows 5-13 (Kraken2_0.0 ... VSEARCH) down to rows 7-15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above the current row 5 (shifts Kraken2_0.0.. down).
$ws.Range("A5:A6").EntireRow.Insert()

# Fill in the new DADA2_Species row (row 5).
$ws.Range("A5").Value = "DADA2_Species"
$ws.Range("B5").Value = 0.03
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0

# Fill in the new DADA2_Taxonomy row (row 6).
$ws.Range("A6").Value = "DADA2_Taxonomy"
$ws.Range("B6").Value = 0.18
$ws.Range("C6").Value = 0.28
$ws.Range("D6").Value = 0.27
$ws.Range("E6").Value = 0.29
$ws.Range("F6").Value = 0.25
